$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1718.5555
$ws.Range("I17").Value = 1165.6666
$ws.Range("J17").Value = 1995
$ws.Range("K17").Value = 3496.9998
$ws.Range("L17").Value = 5985
$ws.Range("M17").Value = -3328.9998
$ws.Range("N17").Value = -6321
$ws.Range("H32").Value = 933.3333
$ws.Range("I32").Value = 900
$ws.Range("K32").Value = 900
$ws.Range("M32").Value = -574
$ws.Range("H40").Value = 1959.8
$ws.Range("I40").Value = 1949.75
$ws.Range("K40").Value = 1949.75
$ws.Range("M40").Value = -1774.75
$ws.Range("H51").Value = 1649.5
$ws.Range("I51").Value = 2800
$ws.Range("J51").Value = 499
$ws.Range("K51").Value = 2800
$ws.Range("L51").Value = 499
$ws.Range("M51").Value = -2316
$ws.Range("N51").Value = -1467
$ws.Range("H62").Value = 7082.4287
$ws.Range("I62").Value = 5901
$ws.Range("K62").Value = 5901
$ws.Range("M62").Value = -5277
$ws.Range("H64").Value = 3033.3333
$ws.Range("I64").Value = 2800
$ws.Range("J64").Value = 3150
$ws.Range("K64").Value = 2800
$ws.Range("L64").Value = 3150
$ws.Range("M64").Value = -2552
$ws.Range("N64").Value = -3646
$ws.Range("H65").Value = 7082.4287
$ws.Range("I65").Value = 5901
$ws.Range("K65").Value = 29505
$ws.Range("M65").Value = -26385
$ws.Range("H67").Value = 3033.3333
$ws.Range("I67").Value = 2800
$ws.Range("J67").Value = 3150
$ws.Range("K67").Value = 2800
$ws.Range("L67").Value = 3150
$ws.Range("M67").Value = -1942
$ws.Range("N67").Value = -4866
$ws.Range("H74").Value = 8250
$ws.Range("H77").Value = 8250
$ws.Range("H80").Value = 297.82352
$ws.Range("J80").Value = 274.6
$ws.Range("L80").Value = 823.8000000000001
$ws.Range("N80").Value = -2819.8
$ws.Range("H83").Value = 297.82352
$ws.Range("J83").Value = 274.6
$ws.Range("L83").Value = 2471.4
$ws.Range("N83").Value = -12455.4
$ws.Range("H116").Value = 5422.857
$ws.Range("I116").Value = 3992
$ws.Range("K116").Value = 3992
$ws.Range("M116").Value = -550

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2725.875
$ws.Range("I2").Value = 2085.6924
$ws.Range("K2").Value = 2085.6924
$ws.Range("M2").Value = -1972.6924
$ws.Range("H32").Value = 16477.795
$ws.Range("I32").Value = 8823.536
$ws.Range("J32").Value = 29872.75
$ws.Range("K32").Value = 8823.536
$ws.Range("L32").Value = 29872.75
$ws.Range("M32").Value = -8536.536
$ws.Range("N32").Value = -30446.75
$ws.Range("H43").Value = 54998
$ws.Range("J43").Value = 54998
$ws.Range("L43").Value = 54998
$ws.Range("N43").Value = -55624
$ws.Range("H116").Value = 2725.875
$ws.Range("I116").Value = 2085.6924
$ws.Range("K116").Value = 2085.6924
$ws.Range("M116").Value = 208.3076000000001
$ws.Range("H122").Value = 913318.8
$ws.Range("I122").Value = 2003301.4
$ws.Range("K122").Value = 6009904.199999999
$ws.Range("M122").Value = -6007454.199999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2725.875
$ws.Range("I3").Value = 2085.6924
$ws.Range("K3").Value = 2085.6924
$ws.Range("M3").Value = -1971.6924
$ws.Range("H82").Value = 66374.5
$ws.Range("J82").Value = 124999
$ws.Range("L82").Value = 124999
$ws.Range("N82").Value = -125765
$ws.Range("H85").Value = 66374.5
$ws.Range("J85").Value = 124999
$ws.Range("L85").Value = 124999
$ws.Range("N85").Value = -127651
$ws.Range("H86").Value = 5311.3335
$ws.Range("I86").Value = 4806
$ws.Range("J86").Value = 5564
$ws.Range("K86").Value = 4806
$ws.Range("L86").Value = 5564
$ws.Range("M86").Value = -3683
$ws.Range("N86").Value = -7810
$ws.Range("H89").Value = 5311.3335
$ws.Range("I89").Value = 4806
$ws.Range("J89").Value = 5564
$ws.Range("K89").Value = 24030
$ws.Range("L89").Value = 27820
$ws.Range("M89").Value = -18414
$ws.Range("N89").Value = -39052

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 583.25
$ws.Range("I22").Value = 167
$ws.Range("K22").Value = 167
$ws.Range("M22").Value = 183
$ws.Range("H31").Value = 5113.636
$ws.Range("J31").Value = 6863.6665
$ws.Range("L31").Value = 6863.6665
$ws.Range("N31").Value = -7453.6665
$ws.Range("H34").Value = 5113.636
$ws.Range("J34").Value = 6863.6665
$ws.Range("L34").Value = 6863.6665
$ws.Range("N34").Value = -7267.6665
$ws.Range("H86").Value = 6640.68
$ws.Range("I86").Value = 3276.5
$ws.Range("J86").Value = 12621.444
$ws.Range("K86").Value = 3276.5
$ws.Range("L86").Value = 12621.444
$ws.Range("M86").Value = -2153.5
$ws.Range("N86").Value = -14867.444
$ws.Range("H89").Value = 6640.68
$ws.Range("I89").Value = 3276.5
$ws.Range("J89").Value = 12621.444
$ws.Range("K89").Value = 16382.5
$ws.Range("L89").Value = 63107.22
$ws.Range("M89").Value = -10766.5
$ws.Range("N89").Value = -74339.22
$ws.Range("H132").Value = 2353.2727
$ws.Range("I132").Value = 2188.8076
$ws.Range("J132").Value = 2964.1428
$ws.Range("K132").Value = 6566.4228
$ws.Range("L132").Value = 8892.428400000001
$ws.Range("M132").Value = -4036.4228
$ws.Range("N132").Value = -13952.4284

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 397.3684
$ws.Range("I5").Value = 397.3684
$ws.Range("K5").Value = 1192.1052
$ws.Range("M5").Value = -1080.1052
$ws.Range("H59").Value = 8532
$ws.Range("I59").Value = 7494
$ws.Range("K59").Value = 22482
$ws.Range("M59").Value = -21942
$ws.Range("H60").Value = 1299.1111
$ws.Range("J60").Value = 1183.1666
$ws.Range("L60").Value = 3549.4998
$ws.Range("N60").Value = -4051.4998
$ws.Range("H113").Value = 4548.3
$ws.Range("I113").Value = 887.5
$ws.Range("J113").Value = 4955.0557
$ws.Range("K113").Value = 2662.5
$ws.Range("L113").Value = 14865.1671
$ws.Range("M113").Value = -492.5
$ws.Range("N113").Value = -19205.1671
$ws.Range("H135").Value = 397.3684
$ws.Range("I135").Value = 397.3684
$ws.Range("K135").Value = 3576.3156
$ws.Range("M135").Value = -1041.3156
$ws.Range("H140").Value = 3851.8125
$ws.Range("I140").Value = 3402.0715
$ws.Range("K140").Value = 10206.2145
$ws.Range("M140").Value = -5026.2145

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 188.38095
$ws.Range("I2").Value = 23.666666
$ws.Range("K2").Value = 23.666666
$ws.Range("M2").Value = 89.33333400000001
$ws.Range("H80").Value = 9747.625
$ws.Range("J80").Value = 9712.286
$ws.Range("L80").Value = 9712.286
$ws.Range("N80").Value = -11708.286
$ws.Range("H83").Value = 9747.625
$ws.Range("J83").Value = 9712.286
$ws.Range("L83").Value = 48561.43
$ws.Range("N83").Value = -58545.43
$ws.Range("H102").Value = 3566.9
$ws.Range("I102").Value = 2225
$ws.Range("J102").Value = 3902.375
$ws.Range("K102").Value = 2225
$ws.Range("L102").Value = 3902.375
$ws.Range("M102").Value = -603
$ws.Range("N102").Value = -7146.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 2857
$ws.Range("I46").Value = 2000
$ws.Range("J46").Value = 2999.8333
$ws.Range("K46").Value = 2000
$ws.Range("L46").Value = 2999.8333
$ws.Range("M46").Value = -1812
$ws.Range("N46").Value = -3375.8333
$ws.Range("H80").Value = 37000
$ws.Range("J80").Value = 37000
$ws.Range("L80").Value = 37000
$ws.Range("N80").Value = -39246
$ws.Range("H82").Value = 2490.5
$ws.Range("I82").Value = 2736.125
$ws.Range("J82").Value = 1999.25
$ws.Range("K82").Value = 2736.125
$ws.Range("L82").Value = 1999.25
$ws.Range("M82").Value = -2375.125
$ws.Range("N82").Value = -2721.25
$ws.Range("H83").Value = 37000
$ws.Range("J83").Value = 37000
$ws.Range("L83").Value = 111000
$ws.Range("N83").Value = -122232
$ws.Range("H85").Value = 2490.5
$ws.Range("I85").Value = 2736.125
$ws.Range("J85").Value = 1999.25
$ws.Range("K85").Value = 2736.125
$ws.Range("L85").Value = 1999.25
$ws.Range("M85").Value = -1488.125
$ws.Range("N85").Value = -4495.25
$ws.Range("H136").Value = 3615.6667
$ws.Range("I136").Value = 3638.8
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 10916.4
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -8366.400000000001
$ws.Range("N136").Value = -15600

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1763.8334
$ws.Range("I122").Value = 1763.8334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 5291.5002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -2841.5002
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2592.6667
$ws.Range("I132").Value = 613
$ws.Range("K132").Value = 1839
$ws.Range("M132").Value = 691
$ws.Range("H136").Value = 85918.25
$ws.Range("I136").Value = 1546.6666
$ws.Range("J136").Value = 170289.83
$ws.Range("K136").Value = 4639.9998
$ws.Range("L136").Value = 510869.49
$ws.Range("M136").Value = -2089.9998
$ws.Range("N136").Value = -515969.49
